# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# --- OFF sheet: update Home (row 2) target depth data ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 219
$wsOff.Range("C2").Value = 158
$wsOff.Range("D2").Value = 39
$wsOff.Range("E2").Value = 16

# --- DEF sheet: update Home (row 2) target depth data ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 216
$wsDef.Range("C2").Value = 142
$wsDef.Range("D2").Value = 50
$wsDef.Range("E2").Value = 27
$wsDef.Range("F2").Value = 5
